$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '64.795.71'
$ws.Range("E2").Value = '  +0.53%  '

Set-TextValue $ws.Range("D3") '3.374.52'
$ws.Range("E3").Value = '  -1.13%  '

Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.30%  '

Set-TextValue $ws.Range("D5") '557.86'
$ws.Range("E5").Value = '  -1.01%  '

Set-TextValue $ws.Range("D6") '177.08'
$ws.Range("E6").Value = '  +2.22%  '

Set-TextValue $ws.Range("D7") '0.622'
$ws.Range("E7").Value = '  +0.36%  '

Set-TextValue $ws.Range("D8") '3.362.29'
$ws.Range("E8").Value = '  -1.35%  '

Set-TextValue $ws.Range("D9") '0.999'
$ws.Range("E9").Value = '  -0.07%  '

Set-TextValue $ws.Range("D10") '0.167'
$ws.Range("E10").Value = '  +7.98%  '

Set-TextValue $ws.Range("D11") '0.633'
$ws.Range("E11").Value = '  +1.84%  '

Set-TextValue $ws.Range("D12") '55.20'
$ws.Range("E12").Value = '  -3.49%  '

Set-TextValue $ws.Range("D13") '0.0000278'
$ws.Range("E13").Value = '  +2.91%  '

Set-TextValue $ws.Range("D14") '9.14'
$ws.Range("E14").Value = '  +0.91%  '

Set-TextValue $ws.Range("D15") '3.906.68'
$ws.Range("E15").Value = '  -1.30%  '

Set-TextValue $ws.Range("D16") '18.36'
$ws.Range("E16").Value = '  +1.91%  '

$ws.Range("E17").Value = '  -1.48%  '

Set-TextValue $ws.Range("D18") '3.370.40'
$ws.Range("E18").Value = '  -1.41%  '

Set-TextValue $ws.Range("D19") '11.89'
$ws.Range("E19").Value = '  +0.58%  '

Set-TextValue $ws.Range("D20") '64.646.16'
$ws.Range("E20").Value = '  +0.24%  '

Set-TextValue $ws.Range("D21") '0.990'
$ws.Range("E21").Value = '  +0.23%  '

Set-TextValue $ws.Range("D22") '457.99'
$ws.Range("E22").Value = '  +11.92%  '

Set-TextValue $ws.Range("D23") '4.83'
$ws.Range("E23").Value = '  +9.80%  '

Set-TextValue $ws.Range("D24") '4.10'
$ws.Range("E24").Value = '  -1.01%  '

Set-TextValue $ws.Range("D25") '86.13'
$ws.Range("E25").Value = '  +3.63%  '

Set-TextValue $ws.Range("D26") '13.47'
$ws.Range("E26").Value = '  +0.61%  '

Set-TextValue $ws.Range("D27") '10.92'
$ws.Range("E27").Value = '  +1.82%  '

Set-TextValue $ws.Range("D28") '2.86'
$ws.Range("E28").Value = '  +3.29%  '

Set-TextValue $ws.Range("D29") '8.83'
$ws.Range("E29").Value = '  -0.25%  '

Set-TextValue $ws.Range("D30") '30.15'
$ws.Range("E30").Value = '  +1.71%  '

Set-TextValue $ws.Range("D31") '6.75'
$ws.Range("E31").Value = '  +0.75%  '

Set-TextValue $ws.Range("D32") '11.51'
$ws.Range("E32").Value = '  +0.26%  '

Set-TextValue $ws.Range("D33") '582.45'
$ws.Range("E33").Value = '  -1.38%  '

Set-TextValue $ws.Range("D34") '0.109'
$ws.Range("E34").Value = '  +1.04%  '

Set-TextValue $ws.Range("D35") '59.55'
$ws.Range("E35").Value = '  +0.99%  '

$ws.Range("E36").Value = '  +0.18%  '

Set-TextValue $ws.Range("D37") '0.142'
$ws.Range("E37").Value = '  -7.89%  '

Set-TextValue $ws.Range("D38") '36.03'
$ws.Range("E38").Value = '  +0.72%  '

$ws.Range("E41").Value = '  +0.26%  '

Set-TextValue $ws.Range("D42") '3.103.51'
$ws.Range("E42").Value = '  -2.00%  '

Set-TextValue $ws.Range("D43") '1.00'
$ws.Range("E43").Value = '  -0.11%  '

Set-TextValue $ws.Range("D44") '2.85'
$ws.Range("E44").Value = '  -1.78%  '

$ws.Range("E45").Value = '  -0.26%  '

Set-TextValue $ws.Range("D46") '0.0413'
$ws.Range("E46").Value = '  +1.08%  '

Set-TextValue $ws.Range("D47") '3.21'
$ws.Range("E47").Value = '  -0.10%  '

$ws.Range("E48").Value = '  +1.47%  '

$ws.Range("E49").Value = '  -1.93%  '

Set-TextValue $ws.Range("D50") '8.31'
$ws.Range("E50").Value = '  -0.19%  '

Set-TextValue $ws.Range("D51") '135.62'
$ws.Range("E51").Value = '  +0.21%  '

# Row 39/40 swap: Stacks (was row 39) <-> PEPE (was row 40)
Set-TextValue $ws.Range("B39") 'PEPE'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D39") '0.0₃0766'
$ws.Range("E39").Value = '  +2.51%  '

Set-TextValue $ws.Range("B40") 'Stacks'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D40") '3.49'
$ws.Range("E40").Value = '  +1.92%  '

